$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-generated K (strikeout) counts for column G, rows 2-71
# (replaces the old Strike# derived values with true K values)
$kValues = @{
    2 = 2
    3 = 1
    4 = 3
    5 = 3
    6 = 2
    7 = 1
    8 = 1
    9 = 3
    10 = 1
    11 = 3
    12 = 2
    13 = 0
    14 = 3
    15 = 0
    16 = 3
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 2
    22 = 2
    23 = 3
    24 = 2
    25 = 1
    26 = 1
    27 = 2
    28 = 2
    29 = 2
    30 = 2
    31 = 2
    32 = 1
    33 = 2
    34 = 2
    35 = 0
    36 = 2
    37 = 3
    38 = 0
    39 = 2
    40 = 1
    41 = 3
    42 = 1
    43 = 1
    44 = 1
    45 = 3
    46 = 1
    47 = 2
    48 = 2
    49 = 2
    50 = 0
    51 = 2
    52 = 3
    53 = 0
    54 = 3
    55 = 3
    56 = 1
    57 = 2
    58 = 2
    59 = 1
    60 = 2
    61 = 1
    62 = 2
    63 = 3
    64 = 3
    65 = 1
    66 = 2
    67 = 2
    68 = 2
    69 = 2
    70 = 2
    71 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
